$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.568.01"
$ws.Range("E2").Value = "  +5.59%  "

$ws.Range("D3").Value = "2.047.54"
$ws.Range("E3").Value = "  +3.24%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'251.78"
$ws.Range("E5").Value = "  +4.03%  "

$ws.Range("D6").Value = "'0.651"
$ws.Range("E6").Value = "  +2.35%  "

$ws.Range("D7").Value = "'65.15"
$ws.Range("E7").Value = "  +16.05%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").Value = "'0.377"
$ws.Range("E9").Value = "  +5.26%  "

$ws.Range("D10").Value = "'59.16"
$ws.Range("E10").Value = "  -1.69%  "

$ws.Range("D11").Value = "'0.0758"
$ws.Range("E11").Value = "  +4.00%  "

$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").Value = "'0.915"
$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("D14").Value = "'15.15"
$ws.Range("E14").Value = "  +6.65%  "

$ws.Range("D15").Value = "2.346.29"
$ws.Range("E15").Value = "  +3.21%  "

$ws.Range("D16").Value = "'5.59"
$ws.Range("E16").Value = "  +6.96%  "

$ws.Range("D17").Value = "'20.44"
$ws.Range("E17").Value = "  +20.18%  "

$ws.Range("D18").Value = "2.055.85"
$ws.Range("E18").Value = "  +3.39%  "

$ws.Range("D19").Value = "37.491.44"
$ws.Range("E19").Value = "  +5.58%  "

$ws.Range("D20").Value = "'73.59"
$ws.Range("E20").Value = "  +4.76%  "

$ws.Range("D21").Value = "0.0₃0876"
$ws.Range("E21").Value = "  +4.91%  "

$ws.Range("D22").Value = "'5.36"
$ws.Range("E22").Value = "  +6.20%  "

$ws.Range("D23").Value = "'237.54"
$ws.Range("E23").Value = "  +2.24%  "

$ws.Range("E24").Value = "  +16.94%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("D27").Value = "'9.60"
$ws.Range("E27").Value = "  +5.94%  "

$ws.Range("D28").Value = "'160.88"
$ws.Range("E28").Value = "  -1.61%  "

$ws.Range("D29").Value = "'19.94"
$ws.Range("E29").Value = "  +2.65%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'5.24"
$ws.Range("E30").Value = "  +9.32%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.122"
$ws.Range("E31").Value = "  +2.63%  "

$ws.Range("E32").Value = "  +24.01%  "

$ws.Range("D33").Value = "'1.21"
$ws.Range("E33").Value = "  +7.40%  "

$ws.Range("D34").Value = "'4.74"
$ws.Range("E34").Value = "  +11.58%  "

$ws.Range("D35").Value = "'0.0614"
$ws.Range("E35").Value = "  +4.91%  "

$ws.Range("D36").Value = "'2.43"
$ws.Range("E36").Value = "  +5.49%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").Value = "'1.85"
$ws.Range("E38").Value = "  +2.87%  "

$ws.Range("D39").Value = "'5.94"
$ws.Range("E39").Value = "  +21.33%  "

$ws.Range("D40").Value = "'0.102"
$ws.Range("E40").Value = "  +15.22%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.81"
$ws.Range("E41").Value = "  +24.93%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.23"
$ws.Range("E42").Value = "  +4.09%  "

$ws.Range("E43").Value = "  +4.52%  "

$ws.Range("D44").Value = "'0.0219"
$ws.Range("E44").Value = "  +4.81%  "

$ws.Range("D45").Value = "'1.14"
$ws.Range("E45").Value = "  +5.42%  "

$ws.Range("E46").Value = "  +9.25%  "

$ws.Range("D47").Value = "'16.95"
$ws.Range("E47").Value = "  +9.95%  "

$ws.Range("D48").Value = "'95.26"
$ws.Range("E48").Value = "  +5.12%  "

$ws.Range("D49").Value = "1.423.65"
$ws.Range("E49").Value = "  +3.84%  "

$ws.Range("D50").Value = "'2.95"
$ws.Range("E50").Value = "  +2.73%  "

$ws.Range("D51").Value = "'47.51"
$ws.Range("E51").Value = "  +3.83%  "
